# Updated testing results for 1-14-2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Code" column is now explicitly "Function Code" ---
$ws.Range("E1").Value = "Function Code"

# --- Row 2 (k-NN test): fill in RAM used ---
$ws.Range("J2").Value = 1.5

# --- Row 4 (MICE cart m=1 run): trim the note down to the stochastic/logged-events
#     remark and move the "wrong dataset" remark into the Errors? column instead ---
$ws.Range("H4").Value = "1. default method is stochastic, while 'cart' is not 2. # of Logged events: 176, "
$ws.Range("I4").Value = "WRONG DATASET USED (non-computational, discovered after running)"
$ws.Range("I4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 43.2

# --- Row 5 (MICE cart re-run, correct dataset): record no errors + RAM used ---
$ws.Range("I5").Value = "None"
$ws.Range("J5").Value = 3.4

# --- Row 6: brand-new test entry for the parallelized MICE (parlmice) run ---
$ws.Range("C6").Value = "2017-18 NHANES Questionnaire Data "
$ws.Range("C6").WrapText = $true
$ws.Range("D6").Value = "9255 x 91"
$ws.Range("E6").Value = 'combined <- parlmice(combined, m = 5, seed = 2022, method = "cart", n.core = 7, n.imp.core = 2, cluster.seed = 1995)'
$ws.Range("E6").WrapText = $true
$ws.Range("H6").Value = "7 cores used, 2 imputations per core. Cluster.seed and seed help ensure reproducibility"
$ws.Range("H6").WrapText = $true
$ws.Range("J6").Value = 6
$ws.Rows.Item(6).RowHeight = 57.6

# --- Sheet view state: scrolled over a bit and left selection on I10 ---
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I10").Select()
